# AMOS B05 - Agile Programming: "Updated after B02 for B03"
#
# Slide 9 contains a 3-column x 6-row table summarizing method categories.
#  - The table style is switched to a different built-in style.
#  - The explicit 12pt font-size override on the 15 body/footer cells
#    (rows 2-6) is removed so the cells fall back to the inherited
#    (14pt) default used by the rest of the deck's table text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

# Apply the new table style.
$tbl.ApplyStyle("{62EB2A9A-D057-4935-8898-D239301CF016}")

# Remove the explicit 12pt size override on every data cell (rows 2-6,
# all 3 columns) so the text reverts to the table's inherited 14pt size.
for ($r = 2; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Shape.TextFrame.TextRange.Font.Size = 14
    }
}
